# live_trading_results.xlsx update:
# Trade #7 closed at 2026-02-17 07:52:54 - unknown UNKNOWN +0.000%
#
# Updates summary / strategy-status aggregate stats for the MarketMaking
# strategy and appends the newly-closed Trade #7 row to both the
# "All Trades" and "MarketMaking" trade logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet - refresh aggregate P&L / trade-count metrics
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.08   # Current Capital
$summary.Range("B4").Value = 0.08      # Total P&L $
$summary.Range("B5").Value = 0.23      # Total P&L %
$summary.Range("B6").Value = 7         # Total Trades
$summary.Range("B7").Value = 4         # Winning Trades
$summary.Range("B9").Value = 57.14     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.08     # Capital
$status.Range("D4").Value = 7          # Trades
$status.Range("E4").Value = 0.08       # P&L $
$status.Range("F4").Value = 0.08       # P&L %
$status.Range("G4").Value = 57.14      # Win Rate %

# ---------------------------------------------------------------
# Append the newly closed Trade #7 to a trade-log sheet
# ---------------------------------------------------------------
function Add-TradeSevenRow($sheet) {
    $sheet.Range("A8").Value = 7
    # "2026-02-17" looks like a date to Excel's smart-entry parser, so
    # pre-format the cell as Text before assigning it, then clear the
    # number-format override again so the cell keeps the default style
    # (same as every other row in this column) while the stored value
    # stays the literal string.
    $sheet.Range("B8").NumberFormat = "@"
    $sheet.Range("B8").Value = "2026-02-17"
    $sheet.Range("B8").ClearFormats()
    $sheet.Range("C8").Value = "07:52:48"
    $sheet.Range("D8").Value = "MarketMaking"
    $sheet.Range("E8").Value = "DOWN"
    $sheet.Range("F8").Value = 0.27
    $sheet.Range("G8").Value = 0.3
    $sheet.Range("H8").Value = "CLOSED"
    $sheet.Range("I8").Value = 11.1111
    $sheet.Range("J8").Value = 0.03
    $sheet.Range("K8").Value = 100.08
    $sheet.Range("L8").Value = 0
    $sheet.Range("M8").Value = 0
    $sheet.Range("N8").Value = 0.6
    $sheet.Range("O8").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P8").Value = "early_exit"
    $sheet.Range("Q8").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeSevenRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeSevenRow $marketMaking
